$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2 through 99
# from 45175 (2023-09-06) to 45183 (2023-09-14)
$ws.Range("C2:C99").Value = 45183
